$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Recursive time (column B) and Iterate time (column C) values
# for disk counts 5, 10, 15, 20 (rows 2-5). The embedded chart's series
# are bound (by formula) to these same ranges, so it keeps tracking the
# new numbers.
$ws.Range("B2").Value = 0.0104888
$ws.Range("C2").Value = 0.0099127
$ws.Range("B3").Value = 1.0106
$ws.Range("C3").Value = 0.972498
$ws.Range("B4").Value = 7.6023
$ws.Range("C4").Value = 3.83543
$ws.Range("B5").Value = 115.516
$ws.Range("C5").Value = 111.271

# Move the active selection to C5, matching the refreshed workbook state.
$ws.Range("C5").Select()
